$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to text format so numeric-looking strings are not
# auto-converted to numbers when the new values are assigned.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '26.090.34'
$ws.Range("E2").Value = '  -1.65%  '
$ws.Range("D3").Value = '1.665.97'
$ws.Range("E3").Value = '  -1.13%  '
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '216.53'
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("D6").Value = '0.5110'
$ws.Range("E6").Value = '  +2.76%  '
$ws.Range("D7").Value = '1.005'
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '0.2634'
$ws.Range("E8").Value = '  +0.94%  '
$ws.Range("D9").Value = '0.06412'
$ws.Range("E9").Value = '  +4.44%  '
$ws.Range("D10").Value = '21.53'
$ws.Range("E10").Value = '  -0.51%  '
$ws.Range("D11").Value = '0.07421'
$ws.Range("E11").Value = '  +2.04%  '
$ws.Range("D12").Value = '1.671.59'
$ws.Range("E12").Value = '  -0.29%  '
$ws.Range("D13").Value = '4.507'
$ws.Range("E13").Value = '  +1.95%  '
$ws.Range("D14").Value = '0.5799'
$ws.Range("E14").Value = '  +1.24%  '
$ws.Range("D15").Value = '0.000008546'
$ws.Range("E15").Value = '  +3.03%  '
$ws.Range("D16").Value = '64.25'
$ws.Range("E16").Value = '  -0.63%  '
$ws.Range("D17").Value = '26.134.25'
$ws.Range("E17").Value = '  -1.59%  '
$ws.Range("D18").Value = '4.926'
$ws.Range("E18").Value = '  -1.57%  '
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("E20").Value = '  +0.96%  '
$ws.Range("D21").Value = '189.24'
$ws.Range("E21").Value = '  +3.63%  '
$ws.Range("D22").Value = '6.207'
$ws.Range("E22").Value = '  +0.68%  '
$ws.Range("D23").Value = '1.007'
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").Value = '145.29'
$ws.Range("E24").Value = '  +0.54%  '
$ws.Range("D25").Value = '7.609'
$ws.Range("E25").Value = '  +0.55%  '
$ws.Range("D26").Value = '0.1200'
$ws.Range("D27").Value = '15.59'
$ws.Range("E27").Value = '  +2.02%  '
$ws.Range("D28").Value = '0.06367'
$ws.Range("E28").Value = '  +13.97%  '
$ws.Range("E29").Value = '  -1.67%  '
$ws.Range("D30").Value = '1.314'
$ws.Range("E30").Value = '  -0.39%  '
$ws.Range("D31").Value = '3.524'
$ws.Range("E31").Value = '  +1.30%  '
$ws.Range("D32").Value = '3.507'
$ws.Range("E32").Value = '  +1.28%  '
$ws.Range("D33").Value = '1.633'
$ws.Range("E33").Value = '  -0.55%  '
$ws.Range("D34").Value = '1.015'
$ws.Range("E34").Value = '  +0.88%  '
$ws.Range("D35").Value = '0.6084'
$ws.Range("E35").Value = '  +3.63%  '
$ws.Range("D36").Value = '2.361'
$ws.Range("E36").Value = '  -0.54%  '
$ws.Range("D37").Value = '2.647'
$ws.Range("E37").Value = '  +0.33%  '
$ws.Range("D38").Value = '6.152'
$ws.Range("E38").Value = '  +4.09%  '
$ws.Range("D39").Value = '0.01606'
$ws.Range("E39").Value = '  +1.52%  '
$ws.Range("D40").Value = '1.075.16'
$ws.Range("E40").Value = '  +0.07%  '
$ws.Range("D41").Value = '0.8596'
$ws.Range("E41").Value = '  +1.21%  '
$ws.Range("E42").Value = '  +0.63%  '
$ws.Range("D43").Value = '101.01'
$ws.Range("E43").Value = '  +2.75%  '
$ws.Range("D44").Value = '1.814.04'
$ws.Range("E44").Value = '  -1.60%  '
$ws.Range("D45").Value = '0.00000000115'
$ws.Range("E45").Value = '  +9.43%  '
$ws.Range("D46").Value = '56.17'
$ws.Range("E46").Value = '  -0.13%  '
$ws.Range("E47").Value = '  +0.23%  '
$ws.Range("D48").Value = '8.052'
$ws.Range("E48").Value = '  -0.06%  '
$ws.Range("D49").Value = '0.05199'
$ws.Range("E49").Value = '  -0.14%  '
$ws.Range("D50").Value = '0.4289'
$ws.Range("E50").Value = '  -0.97%  '
$ws.Range("D51").Value = '5.934'
$ws.Range("E51").Value = '  +6.75%  '

# Restore the original (default) cell style now that the values are set,
# so only the cell contents change and formatting matches the source file.
$dataRange.Style = "Normal"
